$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.115563869476318
$ws.Range("B1").Value = 2.433353424072266
$ws.Range("C1").Value = 9.735664367675781
$ws.Range("D1").Value = 2.205146551132202
$ws.Range("E1").Value = 1.274863719940186
